# Add a new "Save" column (H) to the sheet:
#  - H1 header "Save", formatted like the other header cells (copy G1's style)
#  - H2:H15 filled with 0/1 "Save" flag values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, same look as the rest of row 1 (bold, centered, bordered)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for the new column
$saveValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 0
    6  = 0
    7  = 1
    8  = 1
    9  = 1
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 0
    15 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
